$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Columns.Item(5).ColumnWidth = 29.9777047293527
$ovw.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40
$zh.Range("J2").Value = "2175c337-6dc0-4918-aa9a-89e15c8d7752.0009188c8570ccdc952443e13dc51e2934816f79.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-24 06:59:58"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81bd80e356ab52cf8c73078ba640afbc66dfebf1/e2e/2175c337-6dc0-4918-aa9a-89e15c8d7752.md", "", "", "2175c337-6dc0-4918-aa9a-89e15c8d7752.md")

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
$de.Range("J2").Value = "2175c337-6dc0-4918-aa9a-89e15c8d7752.0009188c8570ccdc952443e13dc51e2934816f79.de-de.xlf"
$de.Range("K2").Value = "2016-08-24 07:00:27"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81bd80e356ab52cf8c73078ba640afbc66dfebf1/e2e/2175c337-6dc0-4918-aa9a-89e15c8d7752.md", "", "", "2175c337-6dc0-4918-aa9a-89e15c8d7752.md")
